$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "63÷3="
$t.Cell(1, 2).Range.Text = "60÷3="
$t.Cell(1, 3).Range.Text = "76÷4="
$t.Cell(1, 4).Range.Text = "78÷6="
$t.Cell(1, 5).Range.Text = "21÷5="

$t.Cell(5, 1).Range.Text = "94÷6="
$t.Cell(5, 2).Range.Text = "49÷9="
$t.Cell(5, 3).Range.Text = "82÷4="
$t.Cell(5, 4).Range.Text = "58÷2="
$t.Cell(5, 5).Range.Text = "31÷5="

$t.Cell(9, 1).Range.Text = "86÷7="
$t.Cell(9, 2).Range.Text = "86÷9="
$t.Cell(9, 3).Range.Text = "67÷2="
$t.Cell(9, 4).Range.Text = "83÷8="
$t.Cell(9, 5).Range.Text = "33÷6="

$t.Cell(13, 1).Range.Text = "76÷8="
$t.Cell(13, 2).Range.Text = "16÷7="
$t.Cell(13, 3).Range.Text = "32÷4="
$t.Cell(13, 4).Range.Text = "59÷5="
$t.Cell(13, 5).Range.Text = "64÷2="

$t.Cell(17, 1).Range.Text = "62÷8="
$t.Cell(17, 2).Range.Text = "98÷7="
$t.Cell(17, 3).Range.Text = "34÷9="
$t.Cell(17, 4).Range.Text = "56÷3="
$t.Cell(17, 5).Range.Text = "97÷7="

Write-Host "Done applying replacements"
